$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New daily rows (482-489): 2021-10-26 through 2021-11-02
# Columns: A DATA | B DESCARTADOS | C EM INVESTIGACAO | D CONFIRMADOS
#          E EXAMINADOS | F RECUPERADOS | G ATIVOS | H HOSPITAL
#          I DOMICILIO | J OBITOS
# ---------------------------------------------------------------------------

$data = @(
    @(44495, 16817, 68, 7031, 23916, 6840, 38, 2, 36, 153),
    @(44496, 16835, 82, 7032, 23949, 6840, 39, 2, 37, 153),
    @(44497, 16917, 39, 7039, 23995, 6848, 38, 2, 36, 153),
    @(44498, 16933, 53, 7040, 24026, 6860, 27, 2, 25, 153),
    @(44499, 16991, 32, 7044, 24067, 6864, 27, 2, 25, 153),
    @(44500, 16991, 32, 7044, 24067, 6869, 22, 2, 20, 153),
    @(44501, 17005, 18, 7045, 24068, 6873, 19, 2, 17, 153),
    @(44502, 17010, 17, 7045, 24072, 6873, 19, 2, 17, 153)
)

$startRow = 482
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# ---------------------------------------------------------------------------
# Delta formulas, filled in two shared-formula blocks (482:483 and 484:489)
# matching how the original data was extended day by day.
# ---------------------------------------------------------------------------

$ws.Range("K482:K483").Formula = "=D482-D481"
$ws.Range("L482:L483").Formula = "=F482-F481"
$ws.Range("M482:M483").Formula = "=B482-B481"
$ws.Range("N482:N483").Formula = "=J482-J481"
$ws.Range("O482:O483").Formula = "=G482-G481"
$ws.Range("P482:P483").Formula = "=C482-C481"

$ws.Range("K484:K489").Formula = "=D484-D483"
$ws.Range("L484:L489").Formula = "=F484-F483"
$ws.Range("M484:M489").Formula = "=B484-B483"
$ws.Range("N484:N489").Formula = "=J484-J483"
$ws.Range("O484:O489").Formula = "=G484-G483"
$ws.Range("P484:P489").Formula = "=C484-C483"

# ---------------------------------------------------------------------------
# View bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------------------

$ws.Range("C495").Select()
